$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B59").Value = 245
$ws.Range("D59").Value = 261
$ws.Range("E59").Value = 182
$ws.Range("B60").Value = 185
$ws.Range("D60").Value = 201
$ws.Range("E60").Value = 182
$ws.Range("B63").Value = 274
$ws.Range("E63").Value = 182
$ws.Range("B64").Value = 303
$ws.Range("E64").Value = 182
$ws.Range("B67").Value = 154
$ws.Range("E67").Value = 182
$ws.Range("B68").Value = 123
$ws.Range("E68").Value = 182
$ws.Range("B71").Value = 333
$ws.Range("E71").Value = 182
$ws.Range("B72").Value = 333
$ws.Range("E72").Value = 182
$ws.Range("B73").Value = 391
$ws.Range("E73").Value = 182
$ws.Range("B76").Value = 90
$ws.Range("E76").Value = 182
$ws.Range("B77").Value = 64
$ws.Range("E77").Value = 182
$ws.Range("B78").Value = 33
$ws.Range("E78").Value = 182
$ws.Range("B81").Value = 394
$ws.Range("B82").Value = 34
$ws.Range("D85").Value = 412
$ws.Range("E85").Value = 221
$ws.Range("B86").Value = 34
$ws.Range("E86").Value = 221
$ws.Range("B93").Value = 424
$ws.Range("E93").Value = 182
$ws.Range("B94").Value = 4
$ws.Range("E94").Value = 182
$ws.Range("D123").Value = 20
$ws.Range("E123").Value = 16
$ws.Range("D124").Value = 440
$ws.Range("E124").Value = 16
$ws.Range("B127").Value = 34
$ws.Range("E127").Value = 96
$ws.Range("B128").Value = 394
$ws.Range("E128").Value = 96
$ws.Range("B131").Value = 64
$ws.Range("E131").Value = 56
$ws.Range("B132").Value = 364
$ws.Range("E132").Value = 56
